# Update the dSF column (column F) values on Sheet1 to reflect the
# repulled / recalculated data described in the commit message
# ("repull data, push all data, mean calculation").
#
# Rows 10, 18, 25, 30 are intentionally left untouched (their dSF value
# did not change between the before/after data pulls).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 2
    4  = 2
    5  = -7
    6  = 7
    7  = -5
    8  = -4
    9  = -3
    11 = -3
    12 = -2
    13 = -1
    14 = -6
    15 = -2
    16 = 1
    17 = 5
    19 = -1
    20 = 2
    21 = -5
    22 = 3
    23 = 7
    24 = -1
    26 = 2
    27 = 2
    28 = -1
    29 = 5
    31 = 4
    32 = -2
    33 = 4
    34 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
